$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'288.23"
$ws.Range("E2").Value = "'-4.06%"

$ws.Range("D3").Value = "'31.01"
$ws.Range("E3").Value = "'-3.67%"

$ws.Range("D4").Value = "'4.885"
$ws.Range("E4").Value = "'-1.43%"

$ws.Range("D5").Value = "'0.07109"
$ws.Range("E5").Value = "'-9.76%"

$ws.Range("D6").Value = "'1.797"
$ws.Range("E6").Value = "'-14.04%"

$ws.Range("D7").Value = "'7.657"
$ws.Range("E7").Value = "'-1.77%"

$ws.Range("D8").Value = "'3.780"
$ws.Range("E8").Value = "'-1.50%"

$ws.Range("D9").Value = "'0.8976"
$ws.Range("E9").Value = "'-3.10%"

$ws.Range("D10").Value = "'0.1645"
$ws.Range("E10").Value = "'-5.31%"

$ws.Range("D11").Value = "'0.07529"
$ws.Range("E11").Value = "'-5.33%"

$ws.Range("D12").Value = "'0.08027"
$ws.Range("E12").Value = "'-7.08%"

$ws.Range("D13").Value = "'0.02992"
$ws.Range("E13").Value = "'-3.62%"

$ws.Range("D14").Value = "'0.09988"
$ws.Range("E14").Value = "'-0.30%"

$ws.Range("D15").Value = "'0.001495"
$ws.Range("E15").Value = "'-2.17%"

$ws.Range("D16").Value = "'0.005690"
$ws.Range("E16").Value = "'-1.79%"

$ws.Range("E18").Value = "'0.31%"

$ws.Range("E19").Value = "'-6.04%"

$ws.Range("E20").Value = "'-0.33%"

$ws.Range("D21").Value = "'0.1299"
$ws.Range("E21").Value = "'-0.84%"

$ws.Range("D22").Value = "'4.268"
$ws.Range("E22").Value = "'-0.81%"

$ws.Range("E23").Value = "'12.18%"

$ws.Range("D24").Value = "'0.04483"
$ws.Range("E24").Value = "'-2.82%"

$ws.Range("D25").Value = "'0.001213"
$ws.Range("E25").Value = "'-1.89%"

$ws.Range("D26").Value = "'0.004647"
$ws.Range("E26").Value = "'4.93%"

$ws.Range("E27").Value = "'0.05%"

$ws.Range("D39").Value = "'0.01621"
$ws.Range("E39").Value = "'-5.09%"

$ws.Range("D40").Value = "'0.04347"
$ws.Range("E40").Value = "'-8.29%"

$ws.Range("D41").Value = "'0.007373"
$ws.Range("E41").Value = "'-0.64%"

$ws.Range("D42").Value = "'0.1304"
$ws.Range("E42").Value = "'-3.64%"

$ws.Range("E43").Value = "'-14.77%"

$ws.Range("D44").Value = "'0.01032"
$ws.Range("E44").Value = "'-8.65%"

$ws.Range("D45").Value = "'0.00005880"
$ws.Range("E45").Value = "'-2.38%"

$ws.Range("E46").Value = "'0.04%"

$ws.Range("D47").Value = "'2.219"
$ws.Range("E47").Value = "'170.52%"

$ws.Range("D49").Value = "'0.00002105"
$ws.Range("E49").Value = "'0.04%"

$ws.Range("D50").Value = "'0.0002005"
$ws.Range("E50").Value = "'0.04%"
